# Update "想去人数" (interested-people count) figures across sheets,
# reflecting refreshed scrape numbers for the same events.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 2436
$ws1.Range("F10").Value = 67
$ws1.Range("F18").Value = 9440
$ws1.Range("F21").Value = 7397
$ws1.Range("F22").Value = 11944
$ws1.Range("F28").Value = 2685
$ws1.Range("F32").Value = 951
$ws1.Range("F36").Value = 4547
$ws1.Range("F41").Value = 557

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F19").Value = 7

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F12").Value = 2436
$ws4.Range("F21").Value = 9440
$ws4.Range("F23").Value = 7397
$ws4.Range("F24").Value = 11944
$ws4.Range("F32").Value = 2685
$ws4.Range("F39").Value = 4547
$ws4.Range("F41").Value = 7
$ws4.Range("F46").Value = 557
